{"js": "// The document has a blank signature line that reads:\n//   \"C.\"  followed by a long run of spaces (a fill-in-the-blank for a name).\n// This fills in the name \"Ricardo Lara Col\u00f3n\" right after \"C. \", keeping the\n// remaining trailing spaces intact, e.g.:\n//   \"C.                                      \"\n//   -> \"C. Ricardo Lara Col\u00f3n                                     \"\n\nconst body = context.document.body;\nconst nameToInsert = \"Ricardo Lara Col\u00f3n\";\n\n// The original run's text is \"C.\" followed by 38 spaces (40 chars total).\n// Search for it precisely (it is unique in the document) instead of relying\n// on a short, ambiguous needle.\nconst blank = \"C.\" + \" \".repeat(38);\nconst hits = body.search(blank, { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nlet target = hits.items.length > 0 ? hits.items[0] : null;\n\n// Fallback: in case the exact run of spaces ever differs, locate the\n// paragraph that is just \"C.\" followed by nothing but whitespace.\nif (!target) {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  const blankParagraph = paragraphs.items.find(\n    (p) => p.text.startsWith(\"C.\") && p.text.slice(2).trim() === \"\" && p.text.length > 2\n  );\n\n  if (!blankParagraph) {\n    throw new Error('Could not find the target \"C.\" placeholder line.');\n  }\n  target = blankParagraph.getRange();\n}\n\n// Narrow down to the literal \"C. \" prefix inside that range so we insert the\n// name right after it (i.e. after the period and the single following\n// space), leaving the rest of the blank space untouched.\nconst prefixHits = target.search(\"C. \", { matchCase: true, matchWholeWord: false });\nprefixHits.load(\"items\");\nawait context.sync();\n\nconst prefix = prefixHits.items[0];\n\n// Collapsed range sitting right after \"C. \".\nconst insertionPoint = prefix.getRange(\"After\");\n\n// Insert the name at that caret position.\ninsertionPoint.insertText(nameToInsert, \"Replace\");\nawait context.sync();\n", "ps1": "# The document has a blank signature line that reads:\n#   \"C.\"  followed by a long run of spaces (a fill-in-the-blank for a name).\n# This fills in the name \"Ricardo Lara Col\u00f3n\" right after \"C. \", keeping the\n# remaining trailing spaces intact, e.g.:\n#   \"C.                                      \"\n#   -> \"C. Ricardo Lara Col\u00f3n                                     \"\n\n$d = $word.ActiveDocument\n$nameToInsert = \"Ricardo Lara Col\u00f3n\"\n\n# Build the 38-space blank run (this PowerShell host doesn't support the\n# string-repeat `*` operator, so build it with a loop instead).\n$blank = \"\"\nfor ($i = 0; $i -lt 38; $i++) { $blank = $blank + \" \" }\n\n$findText = \"C.\" + $blank\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $findText\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$found = $rng.Find.Execute()\n\n$startPos = -1\n\nif ($found) {\n    $startPos = $rng.Start\n} else {\n    # Fallback: in case the exact run of spaces ever differs, scan paragraphs\n    # for one that is just \"C.\" followed by nothing but whitespace.\n    $n = $d.Paragraphs.Count\n    for ($i = 1; $i -le $n; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        if ($t.Length -gt 2 -and $t.Substring(0, 2) -eq \"C.\") {\n            $rest = $t.Substring(2).TrimEnd([char]13, [char]7).Trim()\n            if ($rest -eq \"\") {\n                $startPos = $p.Range.Start\n                break\n            }\n        }\n    }\n}\n\nif ($startPos -lt 0) {\n    throw 'Could not find the target \"C.\" placeholder line.'\n}\n\n# Collapse to the caret right after \"C. \" (period + one space = 3 chars)\n# and insert the name there, leaving the remaining 37 spaces untouched.\n$insertionPoint = $d.Range($startPos + 3, $startPos + 3)\n$insertionPoint.InsertAfter($nameToInsert)\n"}
